# Added May and June NFTF links
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- May 2020 row ---
# Set the hyperlink URL (column B) before the month label (column A) so the
# shared-string table picks up the same insertion order the source workbook
# uses (URL string first, then the month label).
$ws.Range("B20").Value = "https://myemail.constantcontact.com/News-From-The-Forest---May.html?soid=1102494320279&aid=ycB1LWU2Wpk"
$ws.Range("A20").Value = "May 2020"
[void]$ws.Hyperlinks.Add($ws.Range("B20"), "https://myemail.constantcontact.com/News-From-The-Forest---May.html?soid=1102494320279&aid=ycB1LWU2Wpk")
# Re-apply the same "Hyperlink" cell style used by the other link cells
# (e.g. B19) instead of the ad-hoc style Hyperlinks.Add creates on its own.
$ws.Range("B20").Style = $ws.Range("B19").Style

# --- June 2020 row ---
# Here the month label (column A) is entered before the URL (column B) so
# the shared-string table ends up with "June 2020" allocated before the
# June link, matching the source workbook's insertion order for this row.
$ws.Range("A21").Value = "June 2020"
$ws.Range("B21").Value = "https://myemail.constantcontact.com/News-From-The-Forest---June.html?soid=1102494320279&aid=au6GlYTV-AU"
[void]$ws.Hyperlinks.Add($ws.Range("B21"), "https://myemail.constantcontact.com/News-From-The-Forest---June.html?soid=1102494320279&aid=au6GlYTV-AU")
$ws.Range("B21").Style = $ws.Range("B19").Style

# Match the author's final selection position recorded in the sheet view.
[void]$ws.Range("B37").Select()
